$wb = $excel.ActiveWorkbook

# --- "info" sheet: update attendee name in B1 ---
$info = $wb.Worksheets.Item("info")
$info.Range("B1").Value = "김장김치"

# --- "items" sheet: update existing rows and append two new rows ---
$items = $wb.Worksheets.Item("items")

function Set-ItemRow {
    param($ws, $row, $name, $unit, $price, $qty)
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $unit
    $ws.Cells.Item($row, 3).Value = $price
    $ws.Cells.Item($row, 4).Value = $qty
    $ws.Cells.Item($row, 5).Value = $price * $qty
}

Set-ItemRow $items 1 "된장국"     "개" 3000 1
Set-ItemRow $items 2 "부의금가방" "줄" 5000 1
Set-ItemRow $items 3 "종이컵"     "줄" 1500 1
Set-ItemRow $items 4 "에어베게"   "개" 2600 1
Set-ItemRow $items 5 "신라면(컵)" "개" 930  1
Set-ItemRow $items 6 "김밥"       "개" 2000 1
Set-ItemRow $items 7 "삼겹살"     "개" 3000 1
Set-ItemRow $items 8 "대패삼겹살" "개" 2000 1
